# Append the 2025-03-08 price row to every "Solar Prices" sheet.
# Each worksheet has a two-column Date/Price table in A1:B6; we extend
# it by one row (A7/B7) with the new date and that sheet's latest price.
# Row 6's cells are plain text (e.g. "40", "1.19", "5,263"), so we force
# the new cells to Text format before assigning, so Excel doesn't
# auto-convert the date-like / number-like strings.

$wb = $excel.ActiveWorkbook

$newRow = 7
$newDate = "2025-03-08"

$prices = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.295"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,263"
    "Silver Busbar front-side"   = "7,879"
    "Silver finger front-side"   = "7,929"
    "USD_CNY"                    = "7.2647"
}

foreach ($ws in $wb.Worksheets) {
    $price = $prices[$ws.Name]
    if ($null -eq $price) { continue }

    $dateCell = $ws.Cells.Item($newRow, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate

    $priceCell = $ws.Cells.Item($newRow, 2)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
}
